$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1987.5
$ws.Range("H18").Value = 17437.25
$ws.Range("I18").Value = 18571.285
$ws.Range("K18").Value = 18571.285
$ws.Range("M18").Value = -18287.285
$ws.Range("H31").Value = 1241.25
$ws.Range("I31").Value = 1241.25
$ws.Range("K31").Value = 3723.75
$ws.Range("M31").Value = -3493.75
$ws.Range("H33").Value = 130.09091
$ws.Range("I33").Value = 143.44444
$ws.Range("K33").Value = 143.44444
$ws.Range("M33").Value = 85.55556000000001
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H112").Value = 3013.8572
$ws.Range("J112").Value = 3799.6
$ws.Range("L112").Value = 11398.8
$ws.Range("N112").Value = -13614.8
$ws.Range("H113").Value = 4498.6665
$ws.Range("I113").Value = 4498.6665
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4498.6665
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -1244.6665
$ws.Range("H115").Value = 485
$ws.Range("I115").Value = 485
$ws.Range("K115").Value = 1455
$ws.Range("M115").Value = 112
$ws.Range("H116").Value = 5000
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H137").Value = 2071.1428
$ws.Range("I137").Value = 1999.8
$ws.Range("J137").Value = 2249.5
$ws.Range("K137").Value = 5999.4
$ws.Range("L137").Value = 6748.5
$ws.Range("M137").Value = -3449.4
$ws.Range("N137").Value = -11848.5
$ws.Range("H138").Value = 3157.5334
$ws.Range("I138").Value = 843.5
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 2530.5
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = 2609.5
$ws.Range("N138").Value = -22277
$ws.Range("H141").Value = 2920.2632
$ws.Range("I141").Value = 2920.2632
$ws.Range("K141").Value = 8760.7896
$ws.Range("M141").Value = -3580.7896
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2271.0908
$ws.Range("I2").Value = 2498.2856
$ws.Range("J2").Value = 1873.5
$ws.Range("K2").Value = 2498.2856
$ws.Range("L2").Value = 1873.5
$ws.Range("M2").Value = -2385.2856
$ws.Range("N2").Value = -2099.5
$ws.Range("H45").Value = 3850.6
$ws.Range("I45").Value = 2751
$ws.Range("K45").Value = 2751
$ws.Range("M45").Value = -2374
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H116").Value = 2271.0908
$ws.Range("I116").Value = 2498.2856
$ws.Range("J116").Value = 1873.5
$ws.Range("K116").Value = 2498.2856
$ws.Range("L116").Value = 1873.5
$ws.Range("M116").Value = -204.2856000000002
$ws.Range("N116").Value = -6461.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2271.0908
$ws.Range("I3").Value = 2498.2856
$ws.Range("J3").Value = 1873.5
$ws.Range("K3").Value = 2498.2856
$ws.Range("L3").Value = 1873.5
$ws.Range("M3").Value = -2384.2856
$ws.Range("N3").Value = -2101.5
$ws.Range("H19").Value = 1366.3334
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 1999.5
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = 1999.5
$ws.Range("M19").Value = 73
$ws.Range("N19").Value = -2345.5
$ws.Range("H94").Value = 2841
$ws.Range("I94").Value = 2841
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2841
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2390
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 31.142857
$ws.Range("I10").Value = 31.142857
$ws.Range("K10").Value = 31.142857
$ws.Range("M10").Value = 107.857143
$ws.Range("H62").Value = 4730.222
$ws.Range("I62").Value = 4845.5
$ws.Range("K62").Value = 4845.5
$ws.Range("M62").Value = -4221.5
$ws.Range("H65").Value = 4730.222
$ws.Range("I65").Value = 4845.5
$ws.Range("K65").Value = 24227.5
$ws.Range("M65").Value = -21107.5
$ws.Range("H132").Value = 810.7692
$ws.Range("I132").Value = 711.6667
$ws.Range("K132").Value = 2135.0001
$ws.Range("M132").Value = 394.9998999999998
$ws.Range("H134").Value = 3851.7778
$ws.Range("I134").Value = 2952.2856
$ws.Range("K134").Value = 8856.856800000001
$ws.Range("M134").Value = -6321.856800000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30.75
$ws.Range("I2").Value = 14.5
$ws.Range("J2").Value = 63.25
$ws.Range("K2").Value = 87
$ws.Range("L2").Value = 379.5
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = -605.5
$ws.Range("H17").Value = 212.5
$ws.Range("J17").Value = 212.5
$ws.Range("L17").Value = 637.5
$ws.Range("N17").Value = -975.5
$ws.Range("H39").Value = 4333.3335
$ws.Range("J39").Value = 4333.3335
$ws.Range("L39").Value = 13000.0005
$ws.Range("N39").Value = -13588.0005
$ws.Range("H60").Value = 805
$ws.Range("I60").Value = 805
$ws.Range("K60").Value = 2415
$ws.Range("M60").Value = -2164
$ws.Range("H74").Value = 839
$ws.Range("I74").Value = 100
$ws.Range("J74").Value = 1578
$ws.Range("K74").Value = 300
$ws.Range("L74").Value = 4734
$ws.Range("M74").Value = 761
$ws.Range("N74").Value = -6856
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H77").Value = 839
$ws.Range("I77").Value = 100
$ws.Range("J77").Value = 1578
$ws.Range("K77").Value = 900
$ws.Range("L77").Value = 14202
$ws.Range("M77").Value = 4404
$ws.Range("N77").Value = -24810
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H92").Value = 1113.1666
$ws.Range("I92").Value = 1113.1666
$ws.Range("K92").Value = 3339.4998
$ws.Range("M92").Value = -2091.4998
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0
$ws.Range("H121").Value = 1222.2222
$ws.Range("I121").Value = 83.333336
$ws.Range("J121").Value = 1791.6666
$ws.Range("K121").Value = 250.000008
$ws.Range("L121").Value = 5374.9998
$ws.Range("M121").Value = 1059.999992
$ws.Range("N121").Value = -7994.9998
$ws.Range("H131").Value = 3569.3125
$ws.Range("I131").Value = 1610
$ws.Range("J131").Value = 4744.9
$ws.Range("K131").Value = 4830
$ws.Range("L131").Value = 14234.7
$ws.Range("M131").Value = 210
$ws.Range("N131").Value = -24314.7
$ws.Range("H139").Value = 3110
$ws.Range("I139").Value = 3110
$ws.Range("K139").Value = 9330
$ws.Range("M139").Value = -4190
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15120625
$ws.Range("I11").Value = 15137857
$ws.Range("J11").Value = 15000000
$ws.Range("K11").Value = 15137857
$ws.Range("L11").Value = 15000000
$ws.Range("M11").Value = -15137718
$ws.Range("N11").Value = -15000278
$ws.Range("H98").Value = 49999
$ws.Range("J98").Value = 49999
$ws.Range("L98").Value = 49999
$ws.Range("N98").Value = -55989
$ws.Range("H122").Value = 12833.333
$ws.Range("I122").Value = 12833.333
$ws.Range("K122").Value = 38499.999
$ws.Range("M122").Value = -36049.999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4078.3333
$ws.Range("I7").Value = 4078.3333
$ws.Range("K7").Value = 4078.3333
$ws.Range("M7").Value = -3966.3333
$ws.Range("H46").Value = 950.5
$ws.Range("I46").Value = 933.3333
$ws.Range("J46").Value = 1002
$ws.Range("K46").Value = 933.3333
$ws.Range("L46").Value = 1002
$ws.Range("M46").Value = -745.3333
$ws.Range("N46").Value = -1378
$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 5666.6665
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H126").Value = 4078.3333
$ws.Range("I126").Value = 4078.3333
$ws.Range("K126").Value = 12234.9999
$ws.Range("M126").Value = -9764.999899999999
$ws.Range("H132").Value = 5477.4
$ws.Range("I132").Value = 6500
$ws.Range("J132").Value = 3943.5
$ws.Range("K132").Value = 19500
$ws.Range("L132").Value = 11830.5
$ws.Range("M132").Value = -16970
$ws.Range("N132").Value = -16890.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -887
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1000
$ws.Range("L9").ClearContents()
$ws.Range("N9").Value = 0
$ws.Range("M9").Value = -860
$ws.Range("H132").Value = 3330.2942
$ws.Range("I132").Value = 1559.8334
$ws.Range("K132").Value = 4679.5002
$ws.Range("M132").Value = -2149.5002
